$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Previous statement periods (rows 16-26, column E) are reordered so the
# most recent period (2412) is now listed first and the oldest (2402)
# last - i.e. the "Periodo Mora" list is reversed. The mora/value (F) and
# base-salary (G) figures are refreshed from the updated source database.
# ----------------------------------------------------------------------

$periods = @("2412", "2411", "2410", "2409", "2408", "2407", "2406", "2405", "2404", "2403", "2402")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 7).Value = 1619843
}

# Row 16 (now period 2412) and row 26 (now period 2402) keep the F-column
# values that used to belong to the opposite end of the list.
$ws.Cells.Item(16, 6).Value = 64794
$ws.Cells.Item(26, 6).Value = 52000
